$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add a new account-statement row for period 2509 ---
# Insert a new blank row at 18, pushing the existing rows (incl. the signature block) down by one.
$ws.Rows.Item(18).Insert()

# Copy the formatting + values of the existing data row (17, period 2508) into the new row 18,
# then change the period label to 2509. This keeps borders/fills/number-formats identical to the
# other data rows.
$ws.Range("B17:J17").Copy($ws.Range("B18:J18"))
$ws.Range("E18").Value = "2509"

# --- Update the summary figures now that there are 3 periods instead of 2 ---
$ws.Range("F13").Value = 3
$ws.Range("E11").Value = 170820
